$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (unchanged values, just re-affirm) ---
$ws.Range("A1").Value = "content"
$ws.Range("B1").Value = "width"
$ws.Range("C1").Value = "height"
$ws.Range("D1").Value = "bits"
$ws.Range("E1").Value = "address"

# --- Clear the old data block so rows can be rewritten in the new order ---
$ws.Range("A2:E20").Clear()

# --- New row data (content name, width, height) in the new row order ---
$ws.Range("A2").Value = "display_ram_1"
$ws.Range("B2").Value = 640
$ws.Range("C2").Value = 480

$ws.Range("A3").Value = "display_ram_2"
$ws.Range("B3").Value = 640
$ws.Range("C3").Value = 480

$ws.Range("A4").Value = "background1"
$ws.Range("B4").Value = 320
$ws.Range("C4").Value = 480

$ws.Range("A5").Value = "background2"
$ws.Range("B5").Value = 320
$ws.Range("C5").Value = 480

$ws.Range("A6").Value = "help"
$ws.Range("B6").Value = 252
$ws.Range("C6").Value = 105

$ws.Range("A7").Value = "crow"
$ws.Range("B7").Value = 72
$ws.Range("C7").Value = 70

$ws.Range("A8").Value = "holyshit"
$ws.Range("B8").Value = 30
$ws.Range("C8").Value = 66

$ws.Range("A9").Value = "loser"
$ws.Range("B9").Value = 163
$ws.Range("C9").Value = 285

$ws.Range("A10").Value = "person_left1"
$ws.Range("B10").Value = 116
$ws.Range("C10").Value = 211

$ws.Range("A11").Value = "person_left2"
$ws.Range("B11").Value = 109
$ws.Range("C11").Value = 209

$ws.Range("A12").Value = "person_left_3"
$ws.Range("B12").Value = 120
$ws.Range("C12").Value = 205

$ws.Range("A13").Value = "person_middle_1"
$ws.Range("B13").Value = 93
$ws.Range("C13").Value = 213

$ws.Range("A14").Value = "person_middle_2"
$ws.Range("B14").Value = 109
$ws.Range("C14").Value = 214

$ws.Range("A15").Value = "person_midlle_3"
$ws.Range("B15").Value = 106
$ws.Range("C15").Value = 218

$ws.Range("A16").Value = "person_right_1"
$ws.Range("B16").Value = 124
$ws.Range("C16").Value = 218

$ws.Range("A17").Value = "person_right_2"
$ws.Range("B17").Value = 125
$ws.Range("C17").Value = 221

$ws.Range("A18").Value = "person_right_3"
$ws.Range("B18").Value = 122
$ws.Range("C18").Value = 216

$ws.Range("A19").Value = "shit"
$ws.Range("B19").Value = 44
$ws.Range("C19").Value = 23

$ws.Range("A20").Value = "start"
$ws.Range("B20").Value = 290
$ws.Range("C20").Value = 121

# --- Column D: bits = width * height ---
# Rows 2:3 share one formula group, row 4 stands alone, rows 5:20 share another group
$ws.Range("D2:D3").Formula = "=B2*C2"
$ws.Range("D4").Formula = "=B4*C4"
$ws.Range("D5:D20").Formula = "=B5*C5"

# --- Column E: address allocation ---
$ws.Range("E2").Value = 0
$ws.Range("E3").Formula = "=E2+D2"
$ws.Range("E4").Formula = "=E3+D3"
$ws.Range("E5").Formula = "=E4+ROUNDUP(D4/2,0)"
$ws.Range("E6:E20").Formula = "=E5+ROUNDUP(D5/2,0)"

# --- Recalculate so cached values are correct ---
$excel.Calculate()

# --- Selection matches the post-edit state ---
$ws.Range("E4:E20").Select()
